$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test Case No changes from 5 -> 1 ---
$ws.Range("A2").Value = 1

# --- Re-shuffle / rewrite the "Given Input" (C) column for rows 2-7 ---
$ws.Range("C2").Value = "using member function"
$ws.Range("C3").Value = "num1=23"
$ws.Range("C4").Value = "num2=21.5"
$ws.Range("C5").Value = "using friend function"
$ws.Range("C6").Value = "num1=12"
$ws.Range("C7").Value = "num2=13.5"

# Old rows 8-10 in column C no longer hold data
$ws.Range("C8:C10").ClearContents()

# --- Rewrite "Expected Output" / "Actual Output" (D & E) columns for rows 2-10 ---
$ws.Range("D2").Value = "using initialization list"
$ws.Range("E2").Value = "using initialization list"

$ws.Range("D3").Value = "5 53.2"
$ws.Range("E3").Value = "5 53.2"

$ws.Range("D4").Value = "10 20.36"
$ws.Range("E4").Value = "10 20.36"

$ws.Range("D5").Value = "using member function"
$ws.Range("E5").Value = "using member function"

$ws.Range("D6").Value = "num1=23"
$ws.Range("E6").Value = "num1=23"

$ws.Range("D7").Value = "num2=21.5"
$ws.Range("E7").Value = "num2=21.5"

$ws.Range("D8").Value = "using friend function"
$ws.Range("E8").Value = "using friend function"

$ws.Range("D9").Value = "num1=12"
$ws.Range("E9").Value = "num1=12"

$ws.Range("D10").Value = "num2=13.5"
$ws.Range("E10").Value = "num2=13.5"

# --- The remaining rows of this sub-table (11-16) are no longer used ---
$ws.Range("D11:E16").ClearContents()

# --- Move the visible selection/scroll position to reflect the trimmed table ---
$ws.Range("E14").Select()
